$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were inserted above the current row 143, pushing the
# existing rows (143..224) down to (145..226).
$ws.Rows.Item(143).Insert()
$ws.Rows.Item(143).Insert()

# New row 143
$ws.Cells.Item(143, 1).Value = 10
$ws.Cells.Item(143, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(143, 3).Value = "La Araucanía"
$ws.Cells.Item(143, 4).Value = 44488
$ws.Cells.Item(143, 5).Value = 9
$ws.Cells.Item(143, 6).Value = 100112037
$ws.Cells.Item(143, 7).Value = "Cebollín"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 90
$ws.Cells.Item(143, 11).Value = 7000
$ws.Cells.Item(143, 12).Value = 8000
$ws.Cells.Item(143, 13).Value = 7556
$ws.Cells.Item(143, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(143, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(143, 16).Value = 630
$ws.Cells.Item(143, 17).Value = 12
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# New row 144
$ws.Cells.Item(144, 1).Value = 10
$ws.Cells.Item(144, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(144, 3).Value = "La Araucanía"
$ws.Cells.Item(144, 4).Value = 44488
$ws.Cells.Item(144, 5).Value = 9
$ws.Cells.Item(144, 6).Value = 100112037
$ws.Cells.Item(144, 7).Value = "Cebollín"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 110
$ws.Cells.Item(144, 11).Value = 5000
$ws.Cells.Item(144, 12).Value = 5000
$ws.Cells.Item(144, 13).Value = 5000
$ws.Cells.Item(144, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(144, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(144, 16).Value = 417
$ws.Cells.Item(144, 17).Value = 12
$ws.Cells.Item(144, 18).Value = "Hortaliza"
